$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Tue Nov 12 17:34:16 EST 2024",
    "Tue Nov 12 17:34:26 EST 2024",
    "Tue Nov 12 17:34:36 EST 2024",
    "Tue Nov 12 17:34:47 EST 2024",
    "Tue Nov 12 17:34:57 EST 2024",
    "Tue Nov 12 17:35:07 EST 2024",
    "Tue Nov 12 17:35:17 EST 2024",
    "Tue Nov 12 17:35:27 EST 2024",
    "Tue Nov 12 17:35:38 EST 2024",
    "Tue Nov 12 17:35:49 EST 2024",
    "Tue Nov 12 17:35:59 EST 2024",
    "Tue Nov 12 17:36:09 EST 2024",
    "Tue Nov 12 17:36:18 EST 2024",
    "Tue Nov 12 17:36:28 EST 2024",
    "Tue Nov 12 17:36:40 EST 2024",
    "Tue Nov 12 17:36:49 EST 2024",
    "Tue Nov 12 17:36:59 EST 2024",
    "Tue Nov 12 17:37:09 EST 2024",
    "Tue Nov 12 17:37:20 EST 2024",
    "Tue Nov 12 17:37:30 EST 2024",
    "Tue Nov 12 17:37:40 EST 2024",
    "Tue Nov 12 17:37:50 EST 2024",
    "Tue Nov 12 17:38:00 EST 2024",
    "Tue Nov 12 17:38:09 EST 2024",
    "Tue Nov 12 17:38:19 EST 2024",
    "Tue Nov 12 17:38:31 EST 2024",
    "Tue Nov 12 17:38:42 EST 2024",
    "Tue Nov 12 17:38:52 EST 2024",
    "Tue Nov 12 17:39:02 EST 2024",
    "Tue Nov 12 17:39:11 EST 2024",
    "Tue Nov 12 17:39:21 EST 2024",
    "Tue Nov 12 17:39:31 EST 2024",
    "Tue Nov 12 17:39:42 EST 2024",
    "Tue Nov 12 17:39:52 EST 2024",
    "Tue Nov 12 17:40:02 EST 2024",
    "Tue Nov 12 17:40:12 EST 2024",
    "Tue Nov 12 17:40:22 EST 2024",
    "Tue Nov 12 17:40:32 EST 2024",
    "Tue Nov 12 17:40:45 EST 2024",
    "Tue Nov 12 17:40:55 EST 2024",
    "Tue Nov 12 17:41:05 EST 2024",
    "Tue Nov 12 17:41:15 EST 2024",
    "Tue Nov 12 17:41:25 EST 2024",
    "Tue Nov 12 17:41:35 EST 2024",
    "Tue Nov 12 17:41:46 EST 2024",
    "Tue Nov 12 17:41:56 EST 2024",
    "Tue Nov 12 17:42:06 EST 2024"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $values[$i]
}
